$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.748.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.404.62'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '404.02'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.18'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.589'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.682'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.40%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.74'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.09%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.35'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.71'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.424.99'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '11.63'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.788.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000140'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '82.93'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '310.50'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.14'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.83'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +10.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.48'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.06'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.02'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.76'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.92%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.114'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '42.53'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.32'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.37'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.36'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.21%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.322'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +11.98%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.94'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.73'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.93'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.23'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.099.80'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.09%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.35'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.79'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +21.13%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.17%  '
